$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-7: 45233 -> 45243
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Range("C$r")
    $val = $cell.Value()
    if ($val.ToOADate() -eq 45233) {
        $cell.Value = 45243
    }
}
